$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Date naissance")

# Fix the typo in the title column for the "rouge" row: remove the stray
# space before the slash ("pionnier /victime" -> "pionnier/victime").
$ws.Range("C2").Value = "pionnier/victime"

# Move the selection / active cell to C3, matching the saved view state.
$ws.Activate()
$ws.Range("C3").Select()
